$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 55
$ws.Range("I8").Value = 55
$ws.Range("K8").Value = 165
$ws.Range("M8").Value = -26

$ws.Range("H28").Value = 1469.0625
$ws.Range("I28").Value = 1508.25
$ws.Range("J28").Value = 1351.5
$ws.Range("K28").Value = 1508.25
$ws.Range("L28").Value = 1351.5
$ws.Range("M28").Value = -1023.25
$ws.Range("N28").Value = -2321.5

$ws.Range("H64").Value = 3701.9185
$ws.Range("I64").Value = 3547.926
$ws.Range("J64").Value = 3890.9092
$ws.Range("K64").Value = 3547.926
$ws.Range("L64").Value = 3890.9092
$ws.Range("M64").Value = -3299.926
$ws.Range("N64").Value = -4386.9092

$ws.Range("H67").Value = 3701.9185
$ws.Range("I67").Value = 3547.926
$ws.Range("J67").Value = 3890.9092
$ws.Range("K67").Value = 3547.926
$ws.Range("L67").Value = 3890.9092
$ws.Range("M67").Value = -2689.926
$ws.Range("N67").Value = -5606.9092

$ws.Range("H107").Value = 4796.607
$ws.Range("I107").Value = 1359.4117
$ws.Range("J107").Value = 10108.637
$ws.Range("K107").Value = 1359.4117
$ws.Range("L107").Value = 10108.637
$ws.Range("M107").Value = 560.5882999999999
$ws.Range("N107").Value = -13948.637

$ws.Range("H111").Value = 38463776
$ws.Range("I111").Value = 100004190
$ws.Range("J111").Value = 1017.125
$ws.Range("K111").Value = 300012570
$ws.Range("L111").Value = 3051.375
$ws.Range("M111").Value = -300009503
$ws.Range("N111").Value = -9185.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3475
$ws.Range("I86").Value = 2633.3333
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 2633.3333
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -1510.3333
$ws.Range("N86").Value = -8246

$ws.Range("H89").Value = 3475
$ws.Range("I89").Value = 2633.3333
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 13166.6665
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -7550.666499999999
$ws.Range("N89").Value = -41232

$ws.Range("H105").Value = 2942.8765
$ws.Range("I105").Value = 1660
$ws.Range("J105").Value = 2987.628
$ws.Range("K105").Value = 1660
$ws.Range("L105").Value = 2987.628
$ws.Range("M105").Value = 87
$ws.Range("N105").Value = -6481.628000000001

$ws.Range("H107").Value = 2119.2856
$ws.Range("I107").Value = 2107.9412
$ws.Range("J107").Value = 2167.5
$ws.Range("K107").Value = 2107.9412
$ws.Range("L107").Value = 2167.5
$ws.Range("M107").Value = -187.9412000000002
$ws.Range("N107").Value = -6007.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 16666.666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 16666.666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 16666.666
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -17006.666

$ws.Range("H62").Value = 10000000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 10000000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H129").Value = 25999.6
$ws.Range("I129").Value = 10000
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 10000
$ws.Range("L129").Value = 49999
$ws.Range("M129").Value = -5000
$ws.Range("N129").Value = -59999

$ws.Range("H130").Value = 43695
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 43695
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 43695
$ws.Range("N130").Value = -53735

$ws.Range("H131").Value = 29400
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 29400
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 29400
$ws.Range("N131").Value = -39480

$ws.Range("H132").Value = 1302.6154
$ws.Range("I132").Value = 1107.3704
$ws.Range("J132").Value = 1741.9166
$ws.Range("K132").Value = 3322.1112
$ws.Range("L132").Value = 5225.7498
$ws.Range("M132").Value = -792.1112000000003
$ws.Range("N132").Value = -10285.7498

$ws.Range("H133").Value = 47325.668
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 47325.668
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 47325.668
$ws.Range("N133").Value = -52385.668

$ws.Range("H134").Value = 2665.5518
$ws.Range("I134").Value = 1079.1052
$ws.Range("J134").Value = 5679.8
$ws.Range("K134").Value = 3237.3156
$ws.Range("L134").Value = 17039.4
$ws.Range("M134").Value = -702.3155999999999
$ws.Range("N134").Value = -22109.4

$ws.Range("H135").Value = 39666.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39666.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39666.332
$ws.Range("N135").Value = -49806.332

$ws.Range("H137").Value = 30884.137
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 30884.137
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 30884.137
$ws.Range("N137").Value = -41084.137

$ws.Range("H138").Value = 39343.25
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 39343.25
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 39343.25
$ws.Range("N138").Value = -49623.25

$ws.Range("H139").Value = 46500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 46500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 46500
$ws.Range("N139").Value = -56780

$ws.Range("H140").Value = 54841.53
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54841.53
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54841.53
$ws.Range("N140").Value = -65201.53

$ws.Range("H141").Value = 32530
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 32530
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 32530
$ws.Range("N141").Value = -42890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 86.25
$ws.Range("I6").Value = 86.25
$ws.Range("K6").Value = 258.75
$ws.Range("M6").Value = -145.75

$ws.Range("H137").Value = 3488.182
$ws.Range("I137").Value = 1741.5385
$ws.Range("J137").Value = 6011.1113
$ws.Range("K137").Value = 5224.6155
$ws.Range("L137").Value = 18033.3339
$ws.Range("M137").Value = -124.6154999999999
$ws.Range("N137").Value = -28233.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7000.6665
$ws.Range("I5").Value = 5501
$ws.Range("K5").Value = 5501
$ws.Range("M5").Value = -5389

$ws.Range("H17").Value = 13440.25
$ws.Range("I17").Value = 222
$ws.Range("J17").Value = 15328.571
$ws.Range("K17").Value = 222
$ws.Range("L17").Value = 15328.571
$ws.Range("M17").Value = -54
$ws.Range("N17").Value = -15664.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2504950
$ws.Range("I7").Value = 5000500
$ws.Range("J7").Value = 9400
$ws.Range("K7").Value = 5000500
$ws.Range("L7").Value = 9400
$ws.Range("M7").Value = -5000387
$ws.Range("N7").Value = -9626

$ws.Range("H96").Value = 1858.75
$ws.Range("J96").Value = 1910
$ws.Range("L96").Value = 1910
$ws.Range("N96").Value = -4656

$ws.Range("H136").Value = 13546.536
$ws.Range("I136").Value = 21821.834
$ws.Range("J136").Value = 1863.7646
$ws.Range("K136").Value = 65465.50199999999
$ws.Range("L136").Value = 5591.293799999999
$ws.Range("M136").Value = -62915.50199999999
$ws.Range("N136").Value = -10691.2938
